$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric-looking price strings to remain text cells
$textCells = @("D5", "D6", "D14", "D16", "D20", "D21", "D22", "D24", "D31", "D32", "D33", "D40", "D41", "D43", "D44", "D45", "D46", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values from the data refresh
$ws.Range('D2').Value = '42.978.48'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '2.301.34'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '309.29'
$ws.Range('E5').Value = '  -2.75%  '
$ws.Range('D6').Value = '104.53'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('E7').Value = '  -1.86%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  -0.64%  '
$ws.Range('E10').Value = '  -1.53%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('E12').Value = '  -3.29%  '
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').Value = '0.988'
$ws.Range('E14').Value = '  +1.05%  '
$ws.Range('D15').Value = '2.780.31'
$ws.Range('E15').Value = '  +4.41%  '
$ws.Range('D16').Value = '15.32'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '2.286.24'
$ws.Range('E17').Value = '  -1.23%  '
$ws.Range('D18').Value = '42.577.23'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('E19').Value = '  -4.29%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').Value = '13.66'
$ws.Range('E20').Value = '  +2.48%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0000105'
$ws.Range('E21').Value = '  -1.17%  '
$ws.Range('D22').Value = '73.21'
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('D24').Value = '267.07'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('E25').Value = '  -1.68%  '
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('E27').Value = '  +17.51%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('E30').Value = '  -1.93%  '
$ws.Range('D31').Value = '36.18'
$ws.Range('E31').Value = '  -4.55%  '
$ws.Range('D32').Value = '165.72'
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').Value = '0.0858'
$ws.Range('E33').Value = '  -3.67%  '
$ws.Range('E34').Value = '  +2.62%  '
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('E36').Value = '  -3.71%  '
$ws.Range('E37').Value = '  -1.19%  '
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('E39').Value = '  +1.69%  '
$ws.Range('D40').Value = '3.61'
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('D41').Value = '109.42'
$ws.Range('E41').Value = '  +11.80%  '
$ws.Range('E42').Value = '  -3.70%  '
$ws.Range('D43').Value = '71.12'
$ws.Range('E43').Value = '  +1.22%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = '0.226'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('D46').Value = '12.20'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('D47').Value = '1.733.69'
$ws.Range('E47').Value = '  +6.53%  '
$ws.Range('D48').Value = '110.71'
$ws.Range('E48').Value = '  -5.67%  '
$ws.Range('D49').Value = '76.88'
$ws.Range('E49').Value = '  -6.22%  '
$ws.Range('D50').Value = '8.65'
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('E51').Value = '  -3.24%  '

# Restore default style on the forced-text cells (keep only the format override gone)
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
